# Update D Plate Number list:
#  - A7 "D1159"  -> "D5555"
#  - A2 "D914"   -> "D00914"
# (order matters for shared-string table ordering: update A7 before A2)
# Finally move the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "D5555"
$ws.Range("A2").Value = "D00914"

$ws.Range("A2").Select()
